$d = $word.ActiveDocument

# The first paragraph contains the placeholder ID text, split across two runs:
#   Run 1: "**ID__AFFARS_5312_topic_8__ID**"
#   Run 2: " " (a trailing space-only run)
$p1 = $d.Paragraphs.Item(1)

# Replace the whole paragraph text (both runs) with the new ID text, and drop
# the trailing space run in the process, by using wildcard Find/Replace across
# the paragraph's range.
$p1.Range.Find.Execute("**ID__AFFARS_5312_topic_8__ID** ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_5312_302__ID**", 2)

# Re-fetch the (possibly re-indexed) first paragraph and update its formatting.
$p1 = $d.Paragraphs.Item(1)

# Update the left indent from 120 to 225 twips (Word properties are in points;
# 1 point = 20 twips).
$p1.Format.LeftIndent = 225 / 20.0

# Add a paragraph border box (top/left/bottom/right) with 5pt padding (space),
# without an explicit line style/weight.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
